# sample_data.xlsx — "Add files via upload" commit
#
# Renames:
#   sheet1 "therapy"  -> "antimicrobials"
#   sheet3 "cultures" -> "microbiology"
#
# microbiology (sheet3) content: a microscopy/culture result was added on
# the "sputum" specimen row, and the "no growth" placeholder rows were
# replaced with real specimen types (sputum / BAL) reporting "No growth";
# one date was corrected; a trailing blank header cell (D1) was added
# matching the other header cells' bold style.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # therapy -> antimicrobials
$ws2 = $wb.Worksheets.Item(2)   # inflammatory_markers (unchanged data)
$ws3 = $wb.Worksheets.Item(3)   # cultures -> microbiology

$ws1.Name = "antimicrobials"
$ws3.Name = "microbiology"

# --- microbiology (sheet3) data edits -------------------------------------

# Row 2: first aerobic-bottle/S. aureus row becomes a microscopy result.
$ws3.Cells.Item(5, 2).Value = "sputum"
$ws3.Cells.Item(8, 2).Value = "BAL"
$ws3.Cells.Item(2, 2).Value = "microscopy"
$ws3.Cells.Item(2, 3).Value = "Gram+ in clusters"
$ws3.Cells.Item(5, 3).Value = "No growth"
$ws3.Cells.Item(8, 3).Value = "No growth"

# Row 3 date correction: 2023-02-22 -> 2023-02-03
$ws3.Cells.Item(3, 1).Value = 44960

# New blank header cell D1, matching the bold header style of A1:C1.
$ws3.Cells.Item(1, 4).Font.Name = "Arial"
$ws3.Cells.Item(1, 4).Font.Size = 10
$ws3.Cells.Item(1, 4).Font.Bold = $true

# Column C is now wider to fit "Gram+ in clusters".
$ws3.Columns.Item(3).ColumnWidth = 15.08984375

# --- selections / active sheet --------------------------------------------
# Each sheet keeps its own last-used cell; microbiology ends up the
# front-most (active) tab.
$ws1.Range("F11").Select()
$ws2.Range("F21").Select()
$ws3.Range("E10").Select()
